$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4729.909
$ws.Range("I33").Value = 90.53333000000001
$ws.Range("J33").Value = 14671.429
$ws.Range("K33").Value = 90.53333000000001
$ws.Range("L33").Value = 14671.429
$ws.Range("M33").Value = 138.46667
$ws.Range("N33").Value = -15129.429

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2459.375
$ws.Range("I41").Value = 3733.6667
$ws.Range("J41").Value = 1694.8
$ws.Range("K41").Value = 3733.6667
$ws.Range("L41").Value = 1694.8
$ws.Range("M41").Value = -3293.6667
$ws.Range("N41").Value = -2574.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1309.8125
$ws.Range("I53").Value = 2043.5
$ws.Range("J53").Value = 869.6
$ws.Range("K53").Value = 2043.5
$ws.Range("L53").Value = 869.6
$ws.Range("M53").Value = -1406.5
$ws.Range("N53").Value = -2143.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1160.6471
$ws.Range("I70").Value = 954.3333
$ws.Range("J70").Value = 1392.75
$ws.Range("K70").Value = 2862.9999
$ws.Range("L70").Value = 4178.25
$ws.Range("M70").Value = -2592.9999
$ws.Range("N70").Value = -4718.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1160.6471
$ws.Range("I73").Value = 954.3333
$ws.Range("J73").Value = 1392.75
$ws.Range("K73").Value = 2862.9999
$ws.Range("L73").Value = 4178.25
$ws.Range("M73").Value = -1926.9999
$ws.Range("N73").Value = -6050.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 581.2857
$ws.Range("I92").Value = 544
$ws.Range("J92").Value = 761.5
$ws.Range("K92").Value = 544
$ws.Range("L92").Value = 761.5
$ws.Range("M92").Value = 704
$ws.Range("N92").Value = -3257.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7228
$ws.Range("I141").Value = 11659
$ws.Range("J141").Value = 4063
$ws.Range("K141").Value = 34977
$ws.Range("L141").Value = 12189
$ws.Range("M141").Value = -29797
$ws.Range("N141").Value = -22549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2210.8235
$ws.Range("I2").Value = 840.4
$ws.Range("J2").Value = 4168.5713
$ws.Range("K2").Value = 840.4
$ws.Range("L2").Value = 4168.5713
$ws.Range("M2").Value = -727.4
$ws.Range("N2").Value = -4394.5713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8200.666999999999
$ws.Range("J37").Value = 8200.666999999999
$ws.Range("L37").Value = 8200.666999999999
$ws.Range("N37").Value = -8746.666999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1727.2
$ws.Range("I74").Value = 1095.3334
$ws.Range("J74").Value = 2675
$ws.Range("K74").Value = 1095.3334
$ws.Range("L74").Value = 2675
$ws.Range("M74").Value = -221.3334
$ws.Range("N74").Value = -4423

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1727.2
$ws.Range("I77").Value = 1095.3334
$ws.Range("J77").Value = 2675
$ws.Range("K77").Value = 5476.666999999999
$ws.Range("L77").Value = 13375
$ws.Range("M77").Value = -1108.666999999999
$ws.Range("N77").Value = -22111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2210.8235
$ws.Range("I116").Value = 840.4
$ws.Range("J116").Value = 4168.5713
$ws.Range("K116").Value = 840.4
$ws.Range("L116").Value = 4168.5713
$ws.Range("M116").Value = 1453.6
$ws.Range("N116").Value = -8756.5713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2210.8235
$ws.Range("I3").Value = 840.4
$ws.Range("J3").Value = 4168.5713
$ws.Range("K3").Value = 840.4
$ws.Range("L3").Value = 4168.5713
$ws.Range("M3").Value = -726.4
$ws.Range("N3").Value = -4396.5713

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 17606.834
$ws.Range("J6").Value = 17606.834
$ws.Range("L6").Value = 17606.834
$ws.Range("N6").Value = -17832.834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 29000
$ws.Range("J108").Value = 29000
$ws.Range("L108").Value = 29000
$ws.Range("N108").Value = -36680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 36947
$ws.Range("J114").Value = 36947
$ws.Range("L114").Value = 36947
$ws.Range("N114").Value = -45625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 25900
$ws.Range("J116").Value = 25900
$ws.Range("L116").Value = 25900
$ws.Range("N116").Value = -35078

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 26000
$ws.Range("J124").Value = 26000
$ws.Range("L124").Value = 26000
$ws.Range("N124").Value = -35820

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 28390
$ws.Range("I126").Value = 27500
$ws.Range("J126").Value = 28835
$ws.Range("K126").Value = 27500
$ws.Range("L126").Value = 28835
$ws.Range("M126").Value = -22560
$ws.Range("N126").Value = -38715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 90000
$ws.Range("J130").Value = 90000
$ws.Range("L130").Value = 90000
$ws.Range("N130").Value = -100040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("N132").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 13900.667
$ws.Range("J74").Value = 13900.667
$ws.Range("L74").Value = 13900.667
$ws.Range("N74").Value = -15648.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 13900.667
$ws.Range("J77").Value = 13900.667
$ws.Range("L77").Value = 41702.001
$ws.Range("N77").Value = -50438.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 6540.364
$ws.Range("J94").Value = 6540.364
$ws.Range("L94").Value = 6540.364
$ws.Range("N94").Value = -7442.364

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2228.7666
$ws.Range("I99").Value = 2282.625
$ws.Range("J99").Value = 2013.3334
$ws.Range("K99").Value = 2282.625
$ws.Range("L99").Value = 2013.3334
$ws.Range("M99").Value = -784.625
$ws.Range("N99").Value = -5009.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2538.625
$ws.Range("I105").Value = 2077.25
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2077.25
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -330.25
$ws.Range("N105").Value = -6494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2228.7666
$ws.Range("I126").Value = 2282.625
$ws.Range("J126").Value = 2013.3334
$ws.Range("K126").Value = 6847.875
$ws.Range("L126").Value = 6040.0002
$ws.Range("M126").Value = -4377.875
$ws.Range("N126").Value = -10980.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 629.13
$ws.Range("I113").Value = 469.33334
$ws.Range("J113").Value = 707.8357999999999
$ws.Range("K113").Value = 1408.00002
$ws.Range("L113").Value = 2123.5074
$ws.Range("M113").Value = 761.9999800000001
$ws.Range("N113").Value = -6463.5074

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6750
$ws.Range("I136").Value = 1466.6666
$ws.Range("J136").Value = 9920
$ws.Range("K136").Value = 4399.9998
$ws.Range("L136").Value = 29760
$ws.Range("M136").Value = -1849.9998
$ws.Range("N136").Value = -34860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1318.909
$ws.Range("I136").Value = 1272.5714
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 3817.7142
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -1267.7142
$ws.Range("N136").Value = -9300

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N132").ClearContents()
